$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" field text (02/12/2016 -> 02/02/2017)
#    on every slide master, every slide layout (for every design/master) and
#    the handout master.
# ---------------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                if ($sh.HasTextFrame) {
                    if ($sh.TextFrame.TextRange.Text -eq "02/12/2016") {
                        $sh.TextFrame.TextRange.Text = "02/02/2017"
                    }
                }
            }
        }
    }
}

for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $master = $design.SlideMaster

    Update-DatePlaceholder $master.Shapes

    for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
        $layout = $master.CustomLayouts.Item($l)
        Update-DatePlaceholder $layout.Shapes
    }
}

if ($p.HasHandoutMaster) {
    Update-DatePlaceholder $p.HandoutMaster.Shapes
}

# ---------------------------------------------------------------------------
# 2) Reformulations on slide 1 (title slide) - collapse split runs that spell
#    out the same sentence back into a single run each.
# ---------------------------------------------------------------------------

$slide1 = $p.Slides.Item(1)

# Title shape: "Jeux de plateau" <br> "-En réalité " / "augmentée " / "-"
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleMerge = $titleRange.Characters(17, 23)
$titleMerge.Text = $titleMerge.Text

# Subtitle shape: "Jeux de plateau en réalité " / "augmentée."
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange
$subtitleMerge = $subtitleRange.Characters(1, 37)
$subtitleMerge.Text = $subtitleMerge.Text

# ---------------------------------------------------------------------------
# 3) Reformulations on slide 4 (Objectifs du projet / Backlog)
# ---------------------------------------------------------------------------

$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange

# Paragraph 1: "Développer ... en réalité " / "augmentée."
$devMerge = $contentRange.Characters(1, 92)
$devMerge.Text = $devMerge.Text

# Paragraph 5: "Integrer" stays separate, " une " / "IA affrontant le joueur." merge
$contentRange2 = $contentShape.TextFrame.TextRange
$iaMerge = $contentRange2.Characters(177, 29)
$iaMerge.Text = $iaMerge.Text
